$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list (crypto prices in column D).
# Values are entered with a leading apostrophe so Excel keeps them as
# literal text (matching the original inlineStr cells) instead of
# reinterpreting the numeric-looking strings as Number cells.
$ws.Range("D2").Value = "'281.80"
$ws.Range("D3").Value = "'20.70"
$ws.Range("D4").Value = "'6.222"
$ws.Range("D5").Value = "'0.06152"
$ws.Range("D7").Value = "'6.570"
$ws.Range("D8").Value = "'1.501"
$ws.Range("D9").Value = "'0.8189"
$ws.Range("D10").Value = "'0.01385"
$ws.Range("D11").Value = "'0.1643"
$ws.Range("D12").Value = "'0.08384"
$ws.Range("D13").Value = "'0.03543"
$ws.Range("D14").Value = "'0.03194"
$ws.Range("D15").Value = "'0.09142"
$ws.Range("D16").Value = "'3.719"
$ws.Range("D17").Value = "'0.001643"
$ws.Range("D18").Value = "'0.04721"
$ws.Range("D19").Value = "'0.006557"
$ws.Range("D20").Value = "'0.006160"
$ws.Range("D21").Value = "'0.001070"
$ws.Range("D22").Value = "'0.0001611"
$ws.Range("D23").Value = "'3.768"
$ws.Range("D24").Value = "'2.322"
$ws.Range("D25").Value = "'0.3355"
$ws.Range("D41").Value = "'0.007187"
$ws.Range("D42").Value = "'0.004503"
$ws.Range("D43").Value = "'0.1098"
$ws.Range("D44").Value = "'0.01105"
$ws.Range("D45").Value = "'0.00006514"
$ws.Range("D48").Value = "'0.002963"
$ws.Range("D49").Value = "'0.00001901"
